# Apply the "Added System and Business Exception Email Notification" change set.
#
# Summary of edits (per the target diff):
#  1. Settings sheet gains 4 new rows (9-12) of Name/Value pairs for the new
#     BusinessException_MailSubject / BusinessException_MailBody /
#     SystemException_MailSubject / SystemException_MailBody config entries.
#  2. Those new Value cells (B9:B12) wrap text like the existing B5:B8 cells,
#     and the two long, multi-line bodies (B10, B12) grow their row height.
#  3. The previously-active tab (Assets, index 2) stops being the workbook's
#     active tab; Settings becomes active/selected instead (reflecting that
#     editing happened there), while Assets' remembered selection moves on
#     to A19 and Settings' remembered selection moves to B10.

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")
$assets   = $wb.Worksheets.Item("Assets")

# --- 1 & 2: new Settings rows -------------------------------------------------
# Write the new shared strings in the exact order needed so they land at the
# same shared-string table indices as the authored workbook (name columns
# first for both pairs, then the long value bodies, interleaving Business vs
# System exception entries the same way the source content does).
$settings.Range("A9").Value  = "BusinessException_MailSubject"
$settings.Range("A10").Value = "BusinessException_MailBody"
$settings.Range("B9").Value  = "Business Exception - OrderID:<OrderID>"
$settings.Range("B10").Value = "Hi,`nBusiness Exception ocurred`nOrderID:<OrderID> `nBusiness Exception Message:<ExceptionMessage>`nRegards,`nRPA JD Robot."
$settings.Range("B11").Value = "RPA System Exception - OrderID:<OrderID>"
$settings.Range("A11").Value = "SystemException_MailSubject"
$settings.Range("A12").Value = "SystemException_MailBody"
$settings.Range("B12").Value = "Hi,`nSystem Exception ocurred`nOrderID:<OrderID> `nSystem Exception Message:<ExceptionMessage>`nStack Trace/Error Details:<StackTrace>`nRegards,`nRPA JD Robot."

# Match the wrap-text styling already used by B5:B8 on this sheet.
$settings.Range("B9:B12").WrapText = $true

# The long bodies need taller rows (consistent with the other wrapped rows
# on this sheet); the short one-line values keep the default row height.
$settings.Rows.Item(10).RowHeight = 90
$settings.Rows.Item(12).RowHeight = 105

# --- 3: active tab / selections -----------------------------------------------
# Assets was the active/selected tab before; move its remembered selection to
# A19, leave it no longer the active tab.
$assets.Activate() | Out-Null
$assets.Range("A19").Select() | Out-Null

# Settings becomes the active tab, with B10 as the remembered selection.
$settings.Activate() | Out-Null
$settings.Range("B10").Select() | Out-Null
